$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '58.215.47'
$ws.Range("E2").Value = '  -2.48%  '
$ws.Range("D3").Value = '2.568.54'
$ws.Range("E3").Value = '  -2.94%  '
$r = $ws.Range("D4")
$r.NumberFormat = "@"
$r.Value = '1.00'
$r.Style = "Normal"
$ws.Range("E4").Value = '  -0.03%  '
$r = $ws.Range("D5")
$r.NumberFormat = "@"
$r.Value = '542.12'
$r.Style = "Normal"
$ws.Range("E5").Value = '  +1.09%  '
$r = $ws.Range("D6")
$r.NumberFormat = "@"
$r.Value = '142.27'
$r.Style = "Normal"
$ws.Range("E6").Value = '  -1.93%  '
$ws.Range("E7").Value = '  +0.05%  '
$r = $ws.Range("D8")
$r.NumberFormat = "@"
$r.Value = '0.578'
$r.Style = "Normal"
$ws.Range("E8").Value = '  +1.16%  '
$r = $ws.Range("D9")
$r.NumberFormat = "@"
$r.Value = '6.74'
$r.Style = "Normal"
$ws.Range("E9").Value = '  +0.12%  '
$ws.Range("E10").Value = '  -3.24%  '
$ws.Range("E11").Value = '  +2.98%  '
$ws.Range("E12").Value = '  -2.28%  '
$ws.Range("D13").Value = '3.024.74'
$ws.Range("E13").Value = '  -3.16%  '
$ws.Range("D14").Value = '58.142.76'
$ws.Range("E14").Value = '  -2.46%  '
$r = $ws.Range("D15")
$r.NumberFormat = "@"
$r.Value = '20.49'
$r.Style = "Normal"
$ws.Range("E15").Value = '  -3.45%  '
$ws.Range("D16").Value = '2.575.24'
$ws.Range("E16").Value = '  -1.56%  '
$ws.Range("E17").Value = '  -2.33%  '
$ws.Range("E18").Value = '  +0.93%  '
$r = $ws.Range("D19")
$r.NumberFormat = "@"
$r.Value = '333.11'
$r.Style = "Normal"
$ws.Range("E19").Value = '  -3.44%  '
$r = $ws.Range("D20")
$r.NumberFormat = "@"
$r.Value = '9.97'
$r.Style = "Normal"
$ws.Range("E20").Value = '  -2.51%  '
$r = $ws.Range("D21")
$r.NumberFormat = "@"
$r.Value = '6.10'
$r.Style = "Normal"
$ws.Range("E21").Value = '  -3.95%  '
$r = $ws.Range("D23")
$r.NumberFormat = "@"
$r.Value = '66.25'
$r.Style = "Normal"
$ws.Range("E23").Value = '  -0.53%  '
$ws.Range("E24").Value = '  +0.74%  '
$r = $ws.Range("D25")
$r.NumberFormat = "@"
$r.Value = '0.999'
$r.Style = "Normal"
$ws.Range("E25").Value = '  +0.05%  '
$ws.Range("E26").Value = '  -5.66%  '
$r = $ws.Range("D27")
$r.NumberFormat = "@"
$r.Value = '7.00'
$r.Style = "Normal"
$ws.Range("E27").Value = '  -3.93%  '
$r = $ws.Range("D28")
$r.NumberFormat = "@"
$r.Value = '0.999'
$r.Style = "Normal"
$ws.Range("E28").Value = '  +0.05%  '
$ws.Range("D29").Value = '0.0₃0726'
$ws.Range("E29").Value = '  -2.87%  '
$r = $ws.Range("D30")
$r.NumberFormat = "@"
$r.Value = '1.64'
$r.Style = "Normal"
$ws.Range("E30").Value = '  -0.78%  '
$ws.Range("B31").Value = 'Aptos'
$ws.Range("C31").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$r = $ws.Range("D31")
$r.NumberFormat = "@"
$r.Value = '5.91'
$r.Style = "Normal"
$ws.Range("E31").Value = '  +0.99%  '
$ws.Range("B32").Value = 'Monero'
$ws.Range("C32").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$r = $ws.Range("D32")
$r.NumberFormat = "@"
$r.Value = '154.39'
$r.Style = "Normal"
$ws.Range("E32").Value = '  +2.81%  '
$ws.Range("E33").Value = '  -0.89%  '
$ws.Range("E34").Value = '  -3.78%  '
$ws.Range("B35").Value = 'SuiNetwork'
$ws.Range("C35").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$r = $ws.Range("D35")
$r.NumberFormat = "@"
$r.Value = '0.845'
$r.Style = "Normal"
$ws.Range("E35").Value = '  +2.68%  '
$ws.Range("B36").Value = 'ImmutableX'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$r = $ws.Range("D36")
$r.NumberFormat = "@"
$r.Value = '1.08'
$r.Style = "Normal"
$ws.Range("E36").Value = '  -5.63%  '
$ws.Range("B37").Value = 'Fetch.AI'
$ws.Range("C37").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$r = $ws.Range("D37")
$r.NumberFormat = "@"
$r.Value = '0.810'
$r.Style = "Normal"
$ws.Range("E37").Value = '  -3.29%  '
$ws.Range("B38").Value = 'Stacks'
$ws.Range("C38").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$r = $ws.Range("D38")
$r.NumberFormat = "@"
$r.Value = '1.41'
$r.Style = "Normal"
$ws.Range("E38").Value = '  -3.62%  '
$ws.Range("B39").Value = 'Filecoin'
$ws.Range("C39").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$r = $ws.Range("D39")
$r.NumberFormat = "@"
$r.Value = '3.56'
$r.Style = "Normal"
$ws.Range("E39").Value = '  -0.75%  '
$ws.Range("B40").Value = 'Bittensor'
$ws.Range("C40").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$r = $ws.Range("D40")
$r.NumberFormat = "@"
$r.Value = '277.40'
$r.Style = "Normal"
$ws.Range("E40").Value = '  -6.41%  '
$ws.Range("B41").Value = 'FirstDigitalUSD'
$ws.Range("C41").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$r = $ws.Range("D41")
$r.NumberFormat = "@"
$r.Value = '1.00'
$r.Style = "Normal"
$ws.Range("E41").Value = '  +0.12%  '
$ws.Range("B42").Value = 'WhiteBITCoin'
$ws.Range("C42").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$r = $ws.Range("D42")
$r.NumberFormat = "@"
$r.Value = '10.62'
$r.Style = "Normal"
$ws.Range("E42").Value = '  -1.04%  '
$ws.Range("B43").Value = 'Mantle'
$ws.Range("C43").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$r = $ws.Range("D43")
$r.NumberFormat = "@"
$r.Value = '0.585'
$r.Style = "Normal"
$ws.Range("E43").Value = '  -3.09%  '
$ws.Range("B44").Value = 'Stellar'
$ws.Range("C44").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$r = $ws.Range("D44")
$r.NumberFormat = "@"
$r.Value = '0.0940'
$r.Style = "Normal"
$ws.Range("E44").Value = '  -1.36%  '
$r = $ws.Range("D45")
$r.NumberFormat = "@"
$r.Value = '0.0528'
$r.Style = "Normal"
$ws.Range("E45").Value = '  -2.93%  '
$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$r = $ws.Range("D46")
$r.NumberFormat = "@"
$r.Value = '18.34'
$r.Style = "Normal"
$ws.Range("E46").Value = '  -5.46%  '
$ws.Range("B47").Value = 'VeChain'
$ws.Range("C47").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$r = $ws.Range("D47")
$r.NumberFormat = "@"
$r.Value = '0.0225'
$r.Style = "Normal"
$ws.Range("E47").Value = '  -0.75%  '
$ws.Range("B48").Value = 'Maker'
$ws.Range("C48").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D48").Value = '1.903.26'
$ws.Range("E48").Value = '  -3.36%  '
$ws.Range("B49").Value = 'InjectiveProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$r = $ws.Range("D49")
$r.NumberFormat = "@"
$r.Value = '17.72'
$r.Style = "Normal"
$ws.Range("E49").Value = '  -3.41%  '
$ws.Range("B50").Value = 'RenderToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$r = $ws.Range("D50")
$r.NumberFormat = "@"
$r.Value = '4.36'
$r.Style = "Normal"
$ws.Range("E50").Value = '  -4.05%  '
$ws.Range("B51").Value = 'Aave'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$r = $ws.Range("D51")
$r.NumberFormat = "@"
$r.Value = '108.69'
$r.Style = "Normal"
$ws.Range("E51").Value = '  -1.67%  '
